$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.134.11'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '1.835.65'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.19'
$ws.Range('E5').Value = '  -2.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6650'
$ws.Range('E6').Value = '  -4.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2956'
$ws.Range('E8').Value = '  -3.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07351'
$ws.Range('E9').Value = '  -4.55%  '
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07683'
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('D12').Value = '1.844.33'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.021'
$ws.Range('E13').Value = '  -2.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6755'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.14'
$ws.Range('E15').Value = '  -5.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.206'
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').Value = '29.050.60'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008228'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.12'
$ws.Range('E19').Value = '  -3.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.52'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9999'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.310'
$ws.Range('E22').Value = '  -4.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '161.05'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1421'
$ws.Range('E25').Value = '  -5.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.691'
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.04'
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.505'
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.228'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.101'
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05298'
$ws.Range('E32').Value = '  +3.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.859'
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7479'
$ws.Range('E34').Value = '  -3.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.130'
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.682'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').Value = '1.318.17'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01805'
$ws.Range('E38').Value = '  -3.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.716'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9225'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.985'
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9984'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.57'
$ws.Range('E43').Value = '  -2.22%  '
$ws.Range('D44').Value = '1.984.85'
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5167'
$ws.Range('E46').Value = '  -3.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '63.86'
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.761'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.243'
$ws.Range('E49').Value = '  -5.89%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05933'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('B51').Value = 'XinFinNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07324'
$ws.Range('E51').Value = '  +7.25%  '

Write-Host "applied"
